$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(106, 2.35,  2.682, 2.74),
    @(107, 2.534, 2.843, 2.647),
    @(108, 2.554, 3.109, 3.007)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
